$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 to the new date
$ws.Range("B2").Value = "Lunes 24/06/2024"

# Delete rows 3 through 11 (no longer needed)
$ws.Range("A3:B11").EntireRow.Delete()
